$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 31.749936
$ws.Range("H2").Value = 95.249808
$ws.Range("I2").Value = 0.5302851438878331
$ws.Range("J2").Value = 0.5302851438878331
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 2447.58743210384
$ws.Range("R2").Value = 22028.28688893456
$ws.Range("S2").Value = 0.1274708270707201
$ws.Range("T2").Value = 0.1274708270707201
$ws.Range("G3").Value = 31.749936
$ws.Range("H3").Value = 95.249808
$ws.Range("I3").Value = 0.5302851438878331
$ws.Range("J3").Value = 0.5302851438878331
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 3225.159684210944
$ws.Range("R3").Value = 29026.4371578985
$ws.Range("S3").Value = 0.167966940420239
$ws.Range("T3").Value = 0.167966940420239
$ws.Range("G4").Value = 31.749936
$ws.Range("H4").Value = 95.249808
$ws.Range("I4").Value = 0.5302851438878331
$ws.Range("J4").Value = 0.5302851438878331
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 4509.341471618816
$ws.Range("R4").Value = 40584.07324456935
$ws.Range("S4").Value = 0.2348473763968739
$ws.Range("T4").Value = 0.2348473763968739
$ws.Range("I5").Value = 0.3451699599880819
$ws.Range("J5").Value = 0.3451699599880819
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 1593.168629640727
$ws.Range("R5").Value = 14338.51766676654
$ws.Range("S5").Value = 0.08297253050889719
$ws.Range("T5").Value = 0.08297253050889719
$ws.Range("I6").Value = 0.3451699599880819
$ws.Range("J6").Value = 0.3451699599880819
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1093320127339601
$ws.Range("T6").Value = 0.1093320127339601
$ws.Range("I7").Value = 0.3451699599880819
$ws.Range("J7").Value = 0.3451699599880819
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.1528654167452246
$ws.Range("T7").Value = 0.1528654167452246
$ws.Range("I8").Value = 0.1245448961240849
$ws.Range("J8").Value = 0.1245448961240849
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 574.8502027627378
$ws.Range("R8").Value = 5173.65182486464
$ws.Range("S8").Value = 0.02993830979306506
$ws.Range("T8").Value = 0.02993830979306505
$ws.Range("I9").Value = 0.1245448961240849
$ws.Range("J9").Value = 0.1245448961240849
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.0394493894238605
$ws.Range("T9").Value = 0.0394493894238605
$ws.Range("I10").Value = 0.1245448961240849
$ws.Range("J10").Value = 0.1245448961240849
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.05515719690715934
$ws.Range("T10").Value = 0.05515719690715933
